$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E and G (sum) values per row (2..8), taken directly from the
# regenerated sval dataset so the stored doubles match bit-for-bit
# (G is nominally B+C+D+E but is written verbatim to avoid any
# floating point summation-order drift).
$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    3 = @(3.230985683306322, 1.667794583268128, 3.900430680208489,  0.496779210170732, 9.295990156953671)
    4 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    6 = @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.837881874639075)
    7 = @(3.230985683306322, 1.667794583268128, 3.900430680208489,  0.496779210170732, 9.295990156953671)
    8 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
